$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "plot/img" column header to "plot_img"
$ws.Range("G1").Value() = "plot_img"

# Add a new question (category "ge", difficulty level 3) whose answer options
# are used up to disable that difficulty level once exhausted.
$newRows = @(
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 0,                     "T", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 3,                     "F", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 5,                     "F", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 10,                    "F", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 20,                    "F", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 50,                    "F", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", 80,                    "F", "F", "wykres_interaktywny_fejsbuk.png"),
    @(7, "ge", 3, "Ile jest szanys ze to zrobimy?", "125 i jedna trzecia", "F", "F", "wykres_interaktywny_fejsbuk.png")
)

$r = 30
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value() = $row[0]
    $ws.Cells.Item($r, 2).Value() = $row[1]
    $ws.Cells.Item($r, 3).Value() = $row[2]
    $ws.Cells.Item($r, 4).Value() = $row[3]
    $ws.Cells.Item($r, 5).Value() = $row[4]
    $ws.Cells.Item($r, 6).Value() = $row[5]
    $ws.Cells.Item($r, 7).Value() = $row[6]
    $ws.Cells.Item($r, 8).Value() = $row[7]
    $r = $r + 1
}

# Update the view: scroll so row 15 is at the top and select K32, matching
# the refreshed sheet view after the new rows were appended.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow() = 15
$win.ScrollColumn() = 1
$ws.Range("K32").Select()
